$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 27 (item id 16719)
$ws.Range("H27").Value = 284
$ws.Range("J27").Value = 284
$ws.Range("L27").Value = 852
$ws.Range("N27").Value = -1054
# Row 82 (item id 12623)
$ws.Range("H82").Value = 499
$ws.Range("I82").Value = 499
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1497
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1091
$ws.Range("N82").ClearContents()
# Row 85 (item id 12623)
$ws.Range("H85").Value = 499
$ws.Range("I85").Value = 499
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1497
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -93
$ws.Range("N85").ClearContents()
# Row 92 (item id 19901)
$ws.Range("H92").Value = 759.73334
$ws.Range("I92").Value = 759.73334
$ws.Range("K92").Value = 759.73334
$ws.Range("M92").Value = 488.26666
# Row 97 (item id 19885)
$ws.Range("H97").Value = 1359.8
$ws.Range("J97").Value = 1359.8
$ws.Range("L97").Value = 4079.4
$ws.Range("N97").Value = -5071.4
# Row 104 (item id 24263)
$ws.Range("H104").Value = 798.6667
$ws.Range("I104").Value = 798.6667
$ws.Range("K104").Value = 2396.0001
$ws.Range("M104").Value = -649.0001000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (item id 27713)
$ws.Range("H2").Value = 1697.4642
$ws.Range("I2").Value = 1601.619
$ws.Range("K2").Value = 1601.619
$ws.Range("M2").Value = -1488.619
# Row 32 (item id 44147)
$ws.Range("H32").Value = 157096.6
$ws.Range("I32").Value = 147741.53
$ws.Range("K32").Value = 147741.53
$ws.Range("M32").Value = -147454.53
# Row 45 (item id 27714)
$ws.Range("H45").Value = 14455.462
$ws.Range("I45").Value = 12991.842
$ws.Range("K45").Value = 12991.842
$ws.Range("M45").Value = -12614.842
# Row 74 (item id 44000)
$ws.Range("H74").Value = 1643.0264
$ws.Range("I74").Value = 1445.2188
$ws.Range("K74").Value = 1445.2188
$ws.Range("M74").Value = -571.2188000000001
# Row 77 (item id 44000)
$ws.Range("H77").Value = 1643.0264
$ws.Range("I77").Value = 1445.2188
$ws.Range("K77").Value = 7226.094000000001
$ws.Range("M77").Value = -2858.094000000001
# Row 82 (item id 10687)
$ws.Range("H82").Value = 19900
$ws.Range("J82").Value = 19900
$ws.Range("L82").Value = 19900
$ws.Range("N82").Value = -20622
# Row 85 (item id 10687)
$ws.Range("H85").Value = 19900
$ws.Range("J85").Value = 19900
$ws.Range("L85").Value = 19900
$ws.Range("N85").Value = -22396
# Row 103 (item id 18533)
$ws.Range("H103").Value = 25787
$ws.Range("J103").Value = 25787
$ws.Range("L103").Value = 25787
$ws.Range("N103").Value = -28131
# Row 116 (item id 27713)
$ws.Range("H116").Value = 1697.4642
$ws.Range("I116").Value = 1601.619
$ws.Range("K116").Value = 1601.619
$ws.Range("M116").Value = 692.3810000000001
# Row 122 (item id 36168)
$ws.Range("H122").Value = 68016.664
$ws.Range("I122").Value = 68016.664
$ws.Range("K122").Value = 204049.992
$ws.Range("M122").Value = -201599.992
# Row 132 (item id 43997)
$ws.Range("H132").Value = 2212.5454
$ws.Range("I132").Value = 2059.1333
$ws.Range("K132").Value = 6177.3999
$ws.Range("M132").Value = -3647.3999

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (item id 27713)
$ws.Range("H3").Value = 1697.4642
$ws.Range("I3").Value = 1601.619
$ws.Range("K3").Value = 1601.619
$ws.Range("M3").Value = -1487.619
# Row 26 (item id 19535)
$ws.Range("H26").Value = 8285.166999999999
$ws.Range("I26").Value = 3838.4
$ws.Range("J26").Value = 30519
$ws.Range("K26").Value = 3838.4
$ws.Range("L26").Value = 30519
$ws.Range("M26").Value = -3546.4
$ws.Range("N26").Value = -31103
# Row 99 (item id 19943)
$ws.Range("H99").Value = 1619.5769
$ws.Range("J99").Value = 1247
$ws.Range("L99").Value = 1247
$ws.Range("N99").Value = -4243

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (item id 44023)
$ws.Range("H31").Value = 4883.6
$ws.Range("I31").Value = 3157.25
$ws.Range("K31").Value = 3157.25
$ws.Range("M31").Value = -2862.25
# Row 34 (item id 44023)
$ws.Range("H34").Value = 4883.6
$ws.Range("I34").Value = 3157.25
$ws.Range("K34").Value = 3157.25
$ws.Range("M34").Value = -2955.25
# Row 68 (item id 10611)
$ws.Range("H68").Value = 24997
$ws.Range("J68").Value = 24997
$ws.Range("L68").Value = 24997
$ws.Range("N68").Value = -26495
# Row 71 (item id 10611)
$ws.Range("H71").Value = 24997
$ws.Range("J71").Value = 24997
$ws.Range("L71").Value = 74991
$ws.Range("N71").Value = -82479
# Row 105 (item id 19928)
$ws.Range("H105").Value = 4269.067
$ws.Range("I105").Value = 5686.5
$ws.Range("J105").Value = 3753.6365
$ws.Range("K105").Value = 5686.5
$ws.Range("L105").Value = 3753.6365
$ws.Range("M105").Value = -3939.5
$ws.Range("N105").Value = -7247.636500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 10 (item id 4689)
$ws.Range("H10").Value = 653.4286
$ws.Range("I10").Value = 278.27274
$ws.Range("J10").Value = 1066.1
$ws.Range("K10").Value = 834.81822
$ws.Range("L10").Value = 3198.3
$ws.Range("M10").Value = -695.81822
$ws.Range("N10").Value = -3476.3
# Row 69 (item id 12850)
$ws.Range("H69").Value = 59956.547
$ws.Range("J69").Value = 68791.78999999999
$ws.Range("L69").Value = 206375.37
$ws.Range("N69").Value = -207997.37
# Row 70 (item id 12867)
$ws.Range("H70").Value = 3202.4
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12630
# Row 72 (item id 12850)
$ws.Range("H72").Value = 59956.547
$ws.Range("J72").Value = 68791.78999999999
$ws.Range("L72").Value = 619126.11
$ws.Range("N72").Value = -627238.11
# Row 73 (item id 12867)
$ws.Range("H73").Value = 3202.4
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -14184
# Row 113 (item id 27843)
$ws.Range("H113").Value = 666.6667
$ws.Range("I113").Value = 666.6667
$ws.Range("K113").Value = 2000.0001
$ws.Range("M113").Value = 169.9999
# Row 131 (item id 36060)
$ws.Range("H131").Value = 111802.836
$ws.Range("I131").Value = 84600.664
$ws.Range("K131").Value = 253801.992
$ws.Range("M131").Value = -248761.992

$ws = $wb.Worksheets.Item("GSM")
# Row 32 (item id 27215)
$ws.Range("H32").Value = 26289
$ws.Range("J32").Value = 26289
$ws.Range("L32").Value = 26289
$ws.Range("N32").Value = -26881
# Row 42 (item id 27213)
$ws.Range("H42").Value = 106666.664
$ws.Range("J42").Value = 106666.664
$ws.Range("L42").Value = 106666.664
$ws.Range("N42").Value = -107636.664
# Row 102 (item id 36169)
$ws.Range("H102").Value = 3963.7273
$ws.Range("I102").Value = 3429.0557
$ws.Range("K102").Value = 3429.0557
$ws.Range("M102").Value = -1807.0557
# Row 115 (item id 27213)
$ws.Range("H115").Value = 106666.664
$ws.Range("J115").Value = 106666.664
$ws.Range("L115").Value = 106666.664
$ws.Range("N115").Value = -109016.664

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (item id 12565)
$ws.Range("H82").Value = 672.5333000000001
$ws.Range("I82").Value = 282.33334
$ws.Range("K82").Value = 282.33334
$ws.Range("M82").Value = 78.66665999999998
# Row 85 (item id 12565)
$ws.Range("H85").Value = 672.5333000000001
$ws.Range("I85").Value = 282.33334
$ws.Range("K85").Value = 282.33334
$ws.Range("M85").Value = 965.66666
# Row 122 (item id 36247)
$ws.Range("H122").Value = 9360
$ws.Range("I122").Value = 10749.357
$ws.Range("K122").Value = 32248.071
$ws.Range("M122").Value = -29798.071

$ws = $wb.Worksheets.Item("WVR")
# Row 3 (item id 3309)
$ws.Range("H3").Value = 25000124
$ws.Range("I3").Value = 25000124
$ws.Range("K3").Value = 25000124
$ws.Range("M3").Value = -25000010
# Row 32 (item id 3066)
$ws.Range("H32").Value = 6000
$ws.Range("I32").Value = 6000
$ws.Range("K32").Value = 6000
$ws.Range("M32").Value = -5683
